# Fixed workflow: the first 4 data rows (old rows 2-5, Cutoff = 0..3) are
# dropped from both result sheets ("NBR" and "BAR"). Deleting those rows
# shifts the remaining data (old rows 6-20) up so it now occupies rows 2-16,
# and the Cutoff column (A) is then re-numbered sequentially starting at 0.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Remove the first four data rows (rows 2 through 5).
    $ws.Rows("2:5").Delete()

    # Re-index the Cutoff column (A) sequentially from 0 for the
    # remaining 15 data rows (now in rows 2-16).
    for ($i = 0; $i -le 14; $i++) {
        $ws.Cells.Item($i + 2, 1).Value = $i
    }
}
